$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the D:E value cells are treated as text so that numeric-looking
# strings (e.g. "0.9966", "1.000", "29.60") are not auto-converted to numbers.
$valueRange = $ws.Range("D2:E51")
$valueRange.NumberFormat = "@"

$ws.Range("D2").Value = "24.952.27"
$ws.Range("E2").Value = "  -3.72%  "
$ws.Range("E3").Value = "  -6.08%  "
$ws.Range("D4").Value = "0.9966"
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").Value = "236.96"
$ws.Range("E5").Value = "  -4.11%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "0.4728"
$ws.Range("E7").Value = "  -6.07%  "
$ws.Range("D8").Value = "0.2560"
$ws.Range("E8").Value = "  -6.41%  "
$ws.Range("D9").Value = "0.06010"
$ws.Range("E9").Value = "  -2.84%  "
$ws.Range("D10").Value = "0.07033"
$ws.Range("E10").Value = "  -3.02%  "
$ws.Range("D11").Value = "1.638.90"
$ws.Range("E11").Value = "  -6.37%  "
$ws.Range("D12").Value = "14.82"
$ws.Range("E12").Value = "  -2.15%  "
$ws.Range("D13").Value = "0.6162"
$ws.Range("E13").Value = "  -5.68%  "
$ws.Range("D14").Value = "4.362"
$ws.Range("E14").Value = "  -5.85%  "
$ws.Range("D15").Value = "72.74"
$ws.Range("E15").Value = "  -6.19%  "
$ws.Range("D16").Value = "0.9998"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "0.9981"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "24.958.58"
$ws.Range("E18").Value = "  -3.82%  "
$ws.Range("D19").Value = "0.000006584"
$ws.Range("E19").Value = "  -3.46%  "
$ws.Range("D20").Value = "11.15"
$ws.Range("E20").Value = "  -5.88%  "
$ws.Range("D21").Value = "4.403"
$ws.Range("E21").Value = "  +1.40%  "
$ws.Range("D22").Value = "1.846.27"
$ws.Range("E22").Value = "  -7.22%  "
$ws.Range("D23").Value = "8.609"
$ws.Range("E23").Value = "  -0.73%  "
$ws.Range("D24").Value = "5.275"
$ws.Range("E24").Value = "  -2.45%  "
$ws.Range("D25").Value = "133.47"
$ws.Range("E25").Value = "  -2.40%  "
$ws.Range("E26").Value = "  -2.63%  "
$ws.Range("D27").Value = "1.359"
$ws.Range("E27").Value = "  -9.27%  "
$ws.Range("D28").Value = "102.60"
$ws.Range("E28").Value = "  -2.98%  "
$ws.Range("D29").Value = "1.660"
$ws.Range("E29").Value = "  -6.35%  "
$ws.Range("E30").Value = "  -4.23%  "
$ws.Range("D31").Value = "0.07719"
$ws.Range("E31").Value = "  -6.24%  "
$ws.Range("D32").Value = "3.566"
$ws.Range("E32").Value = "  -1.61%  "
$ws.Range("D33").Value = "0.9986"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").Value = "0.04312"
$ws.Range("E34").Value = "  -7.77%  "
$ws.Range("D35").Value = "2.601"
$ws.Range("E35").Value = "  -1.95%  "
$ws.Range("D36").Value = "0.9212"
$ws.Range("E36").Value = "  -7.47%  "
$ws.Range("D37").Value = "0.5827"
$ws.Range("E37").Value = "  -5.82%  "
$ws.Range("D38").Value = "2.573"
$ws.Range("E38").Value = "  -5.94%  "
$ws.Range("E39").Value = "  -3.34%  "
$ws.Range("D40").Value = "0.9986"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "0.8348"
$ws.Range("E41").Value = "  +10.13%  "
$ws.Range("D42").Value = "1.803"
$ws.Range("E42").Value = "  -5.78%  "
$ws.Range("D43").Value = "97.38"
$ws.Range("E43").Value = "  -2.32%  "
$ws.Range("D44").Value = "0.3714"
$ws.Range("E44").Value = "  -4.34%  "
$ws.Range("D45").Value = "4.739"
$ws.Range("E45").Value = "  -5.21%  "
$ws.Range("D46").Value = "0.1102"
$ws.Range("E46").Value = "  -3.76%  "
$ws.Range("D47").Value = "0.05216"
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("D48").Value = "6.088"
$ws.Range("E48").Value = "  -3.51%  "
$ws.Range("D49").Value = "29.60"
$ws.Range("D50").Value = "0.9987"
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("D51").Value = "0.9971"
$ws.Range("E51").Value = "  -0.65%  "

# Restore default (General) formatting look without reverting the stored
# text values back to numbers.
$valueRange.ClearFormats()

